# Update cryptocurrency price (D) and volume change (E) columns
# Values that look like plain numbers are prefixed with a leading
# apostrophe so Excel stores them as text (matching the workbook's
# inline-string / text convention) instead of auto-converting them
# to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.968.96'
$ws.Range("E2").Value = '  +1.99%  '
$ws.Range("D3").Value = '3.635.12'
$ws.Range("E3").Value = '  +3.70%  '
$ws.Range("D4").Value = '''0.998'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''605.02'
$ws.Range("E5").Value = '  +0.29%  '
$ws.Range("D6").Value = '''200.57'
$ws.Range("E6").Value = '  +2.51%  '
$ws.Range("D7").Value = '''0.629'
$ws.Range("E7").Value = '  +1.38%  '
$ws.Range("D8").Value = '''0.999'
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +9.55%  '
$ws.Range("D10").Value = '''0.649'
$ws.Range("E10").Value = '  +0.74%  '
$ws.Range("D11").Value = '''53.77'
$ws.Range("E11").Value = '  +1.26%  '
$ws.Range("D12").Value = '''0.0000306'
$ws.Range("E12").Value = '  +2.56%  '
$ws.Range("D13").Value = '''9.60'
$ws.Range("D14").Value = '4.203.85'
$ws.Range("E14").Value = '  +3.54%  '
$ws.Range("D15").Value = '''681.39'
$ws.Range("E15").Value = '  +13.83%  '
$ws.Range("D16").Value = '70.971.57'
$ws.Range("E16").Value = '  +1.79%  '
$ws.Range("D17").Value = '''12.92'
$ws.Range("E17").Value = '  +2.40%  '
$ws.Range("D18").Value = '3.624.80'
$ws.Range("E18").Value = '  +2.71%  '
$ws.Range("D19").Value = '''19.07'
$ws.Range("E19").Value = '  +0.74%  '
$ws.Range("E20").Value = '  +0.05%  '
$ws.Range("E21").Value = '  +1.97%  '
$ws.Range("D22").Value = '''18.78'
$ws.Range("E22").Value = '  +4.97%  '
$ws.Range("D23").Value = '''5.39'
$ws.Range("E23").Value = '  +1.72%  '
$ws.Range("D24").Value = '''105.20'
$ws.Range("E24").Value = '  +3.48%  '
$ws.Range("E25").Value = '  +0.31%  '
$ws.Range("D26").Value = '''3.04'
$ws.Range("E26").Value = '  -2.07%  '
$ws.Range("D27").Value = '''10.55'
$ws.Range("E27").Value = '  -1.87%  '
$ws.Range("E28").Value = '  +4.33%  '
$ws.Range("D29").Value = '''34.30'
$ws.Range("E29").Value = '  +4.12%  '
$ws.Range("D30").Value = '''4.63'
$ws.Range("E30").Value = '  +8.48%  '
$ws.Range("D31").Value = '''7.24'
$ws.Range("E31").Value = '  +3.92%  '
$ws.Range("E32").Value = '  -0.51%  '
$ws.Range("D33").Value = '''0.116'
$ws.Range("E33").Value = '  +1.66%  '
$ws.Range("D34").Value = '''63.31'
$ws.Range("E34").Value = '  +0.33%  '
$ws.Range("D35").Value = '0.0₃0866'
$ws.Range("E35").Value = '  +7.00%  '
$ws.Range("D36").Value = '3.922.33'
$ws.Range("E36").Value = '  +4.98%  '
$ws.Range("E37").Value = '  -0.22%  '
$ws.Range("D38").Value = '''519.73'
$ws.Range("E38").Value = '  +5.70%  '
$ws.Range("D39").Value = '''3.02'
$ws.Range("E39").Value = '  -4.56%  '
$ws.Range("D40").Value = '''0.391'
$ws.Range("E40").Value = '  +0.69%  '
$ws.Range("D41").Value = '''3.58'
$ws.Range("E41").Value = '  -1.26%  '
$ws.Range("D42").Value = '''36.54'
$ws.Range("E42").Value = '  +1.47%  '
$ws.Range("D43").Value = '''0.136'
$ws.Range("E43").Value = '  +2.71%  '
$ws.Range("D44").Value = '''0.0458'
$ws.Range("E44").Value = '  +2.14%  '
$ws.Range("D45").Value = '''3.06'
$ws.Range("E45").Value = '  +8.83%  '
$ws.Range("D46").Value = '''3.45'
$ws.Range("E46").Value = '  +6.36%  '
$ws.Range("E47").Value = '  +1.63%  '
$ws.Range("D48").Value = '''8.65'
$ws.Range("E48").Value = '  +3.11%  '
$ws.Range("E49").Value = '  -0.14%  '
$ws.Range("D50").Value = '''0.000248'
$ws.Range("E50").Value = '  +2.00%  '
$ws.Range("D51").Value = '''1.30'
$ws.Range("E51").Value = '  +2.49%  '
